$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5777637958526611
$ws.Range("B1").Value = 0.5673599243164062
$ws.Range("C1").Value = 0.5918037891387939
$ws.Range("D1").Value = 0.7931169867515564
$ws.Range("E1").Value = 0.8122955560684204
